# Update "Forecast Comparison" sheet with the corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week"
#  - drop the zero-padding on the Week labels (W01 -> W1, ... W09 -> W9)
#  - populate the new Week_Start_Date column with the weekly start dates
#  - store is_holiday_week as a boolean instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (C..J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Corrected (non zero-padded) week labels.
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

# Weekly start dates aligned with each week label. Stored as plain text
# (not auto-converted to a date serial) to match the source data.
$weekStartDates = @(
  "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
  "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
  "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
  "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

# Make sure text typed into column B is kept as text, not parsed as a date.
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $weeks[$i]
    $ws.Cells.Item($r, 2).Value = $weekStartDates[$i]
    # is_holiday_week (now column J) becomes a boolean value.
    $current = [int]($ws.Cells.Item($r, 10).Value)
    $ws.Cells.Item($r, 10).Value = [bool]$current
}
